# Fix NA problem: there is no unemployment data for this series before 2003.
# Change the time range covered by the sheet from 2000-2020 to 2003-2020,
# i.e. remove the first three data rows (years 2000, 2001, 2002).

$wb = $excel.ActiveWorkbook

# "BLS Data Series" is the worksheet that holds the Year / Jan..Dec grid
# (this is the sheet stored as xl/worksheets/sheet1.xml).
$data = $wb.Worksheets.Item(1)

# Select rows 2-4 (years 2000, 2001, 2002) and delete them, shifting the
# remaining years (2003-2020) up so the table starts at row 2 again.
$data.Rows("2:4").Select()
$data.Rows("2:4").Delete()

# "Sheet1" is the worksheet with the series metadata block
# (stored as xl/worksheets/sheet2.xml). A couple of its description rows
# grew taller (e.g. the "Labor force status:" / "Type of data:" rows).
$meta = $wb.Worksheets.Item(2)
$meta.Rows(7).RowHeight = 28
$meta.Rows(8).RowHeight = 28
